$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1288.125
$ws.Range("I11").Value = 1288.125
$ws.Range("K11").Value = 1288.125
$ws.Range("M11").Value = -1148.125
$ws.Range("H19").Value = 2050.6924
$ws.Range("J19").Value = 2632.6667
$ws.Range("L19").Value = 2632.6667
$ws.Range("N19").Value = -2982.6667
$ws.Range("H33").Value = 3220.1155
$ws.Range("I33").Value = 3471.4167
$ws.Range("J33").Value = 204.5
$ws.Range("K33").Value = 3471.4167
$ws.Range("L33").Value = 204.5
$ws.Range("M33").Value = -3242.4167
$ws.Range("N33").Value = -662.5
$ws.Range("H53").Value = 4202.222
$ws.Range("I53").Value = 147.46666
$ws.Range("J53").Value = 9270.666999999999
$ws.Range("K53").Value = 147.46666
$ws.Range("L53").Value = 9270.666999999999
$ws.Range("M53").Value = 489.53334
$ws.Range("N53").Value = -10544.667
$ws.Range("H106").Value = 2970.625
$ws.Range("I106").Value = 2174.4
$ws.Range("K106").Value = 2174.4
$ws.Range("M106").Value = -1543.4
$ws.Range("H135").Value = 685.225
$ws.Range("I135").Value = 641.2059
$ws.Range("J135").Value = 934.6667
$ws.Range("K135").Value = 5770.8531
$ws.Range("L135").Value = 8412.0003
$ws.Range("M135").Value = -3235.8531
$ws.Range("N135").Value = -13482.0003
$ws.Range("H137").Value = 53060.285
$ws.Range("I137").Value = 70063.53999999999
$ws.Range("J137").Value = 3939.7778
$ws.Range("K137").Value = 210190.62
$ws.Range("L137").Value = 11819.3334
$ws.Range("M137").Value = -207640.62
$ws.Range("N137").Value = -16919.3334
$ws.Range("H138").Value = 3497.2334
$ws.Range("I138").Value = 1302.0667
$ws.Range("K138").Value = 3906.2001
$ws.Range("M138").Value = 1233.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1743.9
$ws.Range("I61").Value = 1743.9
$ws.Range("K61").Value = 1743.9
$ws.Range("M61").Value = -1531.9
$ws.Range("H74").Value = 73936.5
$ws.Range("I74").Value = 4298.909
$ws.Range("K74").Value = 4298.909
$ws.Range("M74").Value = -3424.909
$ws.Range("H77").Value = 73936.5
$ws.Range("I77").Value = 4298.909
$ws.Range("K77").Value = 21494.545
$ws.Range("M77").Value = -17126.545
$ws.Range("H102").Value = 3791159.5
$ws.Range("I102").Value = 6412376.5
$ws.Range("K102").Value = 6412376.5
$ws.Range("M102").Value = -6410754.5
$ws.Range("H122").Value = 720403.4399999999
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250
$ws.Range("H132").Value = 3306.3872
$ws.Range("I132").Value = 2432.3635
$ws.Range("K132").Value = 7297.0905
$ws.Range("M132").Value = -4767.0905
$ws.Range("H135").Value = 1000000000
$ws.Range("J135").Value = 1000000000
$ws.Range("L135").Value = 1000000000
$ws.Range("N135").Value = -1000010140
$ws.Range("H136").Value = 1743.9
$ws.Range("I136").Value = 1743.9
$ws.Range("K136").Value = 5231.700000000001
$ws.Range("M136").Value = -2681.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 59992
$ws.Range("J109").Value = 59992
$ws.Range("L109").Value = 59992
$ws.Range("N109").Value = -62766
$ws.Range("H110").Value = 32740.2
$ws.Range("J110").Value = 32740.2
$ws.Range("L110").Value = 32740.2
$ws.Range("N110").Value = -40920.2
$ws.Range("H134").Value = 4692.5
$ws.Range("I134").Value = 1846.8572
$ws.Range("K134").Value = 5540.571599999999
$ws.Range("M134").Value = -3005.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14273.648
$ws.Range("I31").Value = 2121.625
$ws.Range("K31").Value = 2121.625
$ws.Range("M31").Value = -1826.625
$ws.Range("H34").Value = 14273.648
$ws.Range("I34").Value = 2121.625
$ws.Range("K34").Value = 2121.625
$ws.Range("M34").Value = -1919.625
$ws.Range("H105").Value = 3441.6
$ws.Range("I105").Value = 3240.818
$ws.Range("K105").Value = 3240.818
$ws.Range("M105").Value = -1493.818
$ws.Range("H107").Value = 1083.4138
$ws.Range("J107").Value = 960.1
$ws.Range("L107").Value = 960.1
$ws.Range("N107").Value = -4800.1
$ws.Range("H122").Value = 3274.8333
$ws.Range("I122").Value = 2987.875
$ws.Range("J122").Value = 3848.75
$ws.Range("K122").Value = 8963.625
$ws.Range("L122").Value = 11546.25
$ws.Range("M122").Value = -6513.625
$ws.Range("N122").Value = -16446.25
$ws.Range("H132").Value = 74585.66
$ws.Range("J132").Value = 446444
$ws.Range("L132").Value = 1339332
$ws.Range("N132").Value = -1344392
$ws.Range("H134").Value = 3034.739
$ws.Range("I134").Value = 1259.091
$ws.Range("J134").Value = 4662.4165
$ws.Range("K134").Value = 3777.273
$ws.Range("L134").Value = 13987.2495
$ws.Range("M134").Value = -1242.273
$ws.Range("N134").Value = -19057.2495
$ws.Range("H135").Value = 98071.10000000001
$ws.Range("J135").Value = 98071.10000000001
$ws.Range("L135").Value = 98071.10000000001
$ws.Range("N135").Value = -108211.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 121.333336
$ws.Range("I2").Value = 149.91667
$ws.Range("J2").Value = 64.166664
$ws.Range("K2").Value = 899.5000200000001
$ws.Range("L2").Value = 384.999984
$ws.Range("M2").Value = -786.5000200000001
$ws.Range("N2").Value = -610.999984
$ws.Range("H12").Value = 59352.465
$ws.Range("I12").Value = 127036.43
$ws.Range("J12").Value = 129
$ws.Range("K12").Value = 381109.29
$ws.Range("L12").Value = 387
$ws.Range("M12").Value = -380936.29
$ws.Range("N12").Value = -733
$ws.Range("H16").Value = 156.42857
$ws.Range("J16").Value = 156.42857
$ws.Range("L16").Value = 469.28571
$ws.Range("N16").Value = -815.28571
$ws.Range("H56").Value = 35722856
$ws.Range("I56").Value = 35722856
$ws.Range("K56").Value = 35722856
$ws.Range("M56").Value = -35722326
$ws.Range("H92").Value = 5900
$ws.Range("I92").Value = 5900
$ws.Range("K92").Value = 17700
$ws.Range("M92").Value = -16452
$ws.Range("H113").Value = 3138.12
$ws.Range("J113").Value = 2319.158
$ws.Range("L113").Value = 6957.474
$ws.Range("N113").Value = -11297.474
$ws.Range("H114").Value = 4351.8
$ws.Range("J114").Value = 5007.75
$ws.Range("L114").Value = 15023.25
$ws.Range("N114").Value = -21531.25
$ws.Range("H131").Value = 11575974
$ws.Range("J131").Value = 12822972
$ws.Range("L131").Value = 38468916
$ws.Range("N131").Value = -38478996
$ws.Range("H132").Value = 2091
$ws.Range("I132").Value = 811.6
$ws.Range("K132").Value = 7304.400000000001
$ws.Range("M132").Value = -4774.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3676.1482
$ws.Range("I132").Value = 3171.3684
$ws.Range("K132").Value = 9514.1052
$ws.Range("M132").Value = -6984.1052
$ws.Range("H134").Value = 28140.625
$ws.Range("J134").Value = 28140.625
$ws.Range("L134").Value = 84421.875
$ws.Range("N134").Value = -89491.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 2553.5
$ws.Range("I17").Value = 2553.5
$ws.Range("K17").Value = 2553.5
$ws.Range("M17").Value = -2383.5
$ws.Range("H62").Value = 12999.5
$ws.Range("J62").Value = 12000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 12999.5
$ws.Range("J65").Value = 12000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42240
$ws.Range("H109").Value = 59989.5
$ws.Range("J109").Value = 59989.5
$ws.Range("L109").Value = 59989.5
$ws.Range("N109").Value = -62763.5
$ws.Range("H133").Value = 398998
$ws.Range("J133").Value = 398998
$ws.Range("L133").Value = 398998
$ws.Range("N133").Value = -404058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 79995.5
$ws.Range("J109").Value = 79995.5
$ws.Range("L109").Value = 79995.5
$ws.Range("N109").Value = -82769.5
$ws.Range("H126").Value = 2185
$ws.Range("I126").Value = 2185
$ws.Range("K126").Value = 6555
$ws.Range("M126").Value = -4085
$ws.Range("H136").Value = 2421.2742
$ws.Range("I136").Value = 1776.9667
$ws.Range("K136").Value = 5330.9001
$ws.Range("M136").Value = -2780.9001

Write-Host "Applied all market price updates"